$p = $ppt.ActivePresentation
$nm = $p.NotesMaster
Write-Output "NotesMaster CustomLayouts Count: $($nm.CustomLayouts.Count)"
